$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, using the same style as other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in Save values for rows 2-8
$saveValues = @(0, 1, 1, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
